# Apply the scraped price (column D) / 1h-volume (column E) updates to the
# "cryptos" worksheet, and (for the BabyDogeCoin insertion) the new Coin/Link
# text in columns B/C for rows 49-51.
#
# Columns D/E hold plain text in the workbook, not numbers (e.g. "1.000",
# "29.249.21", "0.08250", "  -0.20%  "). Excel auto-converts a literal like
# "1.000" typed into .Value to the number 1, so for any new D-column value that
# Excel would otherwise reinterpret as a number we prefix it with a quote
# (the same quote-prefix Excel's UI uses to force text entry) and then reset
# the cell style to "Normal" afterwards so no numeric style is left applied -
# the result is a plain text cell with the exact literal value, like the source.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.249.21'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '1.866.16'
$ws.Range('E3').Value = '  +0.36%  '
$ws.Range('D4').Value = "'" + '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'" + '0.7228'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.19%  '
$ws.Range('D6').Value = "'" + '240.82'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.41%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = "'" + '0.07838'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.83%  '
$ws.Range('D9').Value = "'" + '0.3088'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.53%  '
$ws.Range('D10').Value = "'" + '25.27'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.06%  '
$ws.Range('D11').Value = "'" + '0.08250'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.26%  '
$ws.Range('D12').Value = '1.868.77'
$ws.Range('E12').Value = '  +0.77%  '
$ws.Range('D13').Value = "'" + '0.7220'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.30%  '
$ws.Range('E14').Value = '  +0.59%  '
$ws.Range('D15').Value = "'" + '90.72'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.72%  '
$ws.Range('D16').Value = '29.291.08'
$ws.Range('E16').Value = '  -0.25%  '
$ws.Range('D17').Value = "'" + '5.858'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.48%  '
$ws.Range('D18').Value = "'" + '243.75'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.44%  '
$ws.Range('D19').Value = "'" + '0.000007808'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.60%  '
$ws.Range('D20').Value = "'" + '13.19'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.32%  '
$ws.Range('D21').Value = '2.108.47'
$ws.Range('E21').Value = '  -0.78%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').Value = "'" + '7.971'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.98%  '
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').Value = "'" + '0.1594'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +11.00%  '
$ws.Range('D26').Value = "'" + '161.88'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.31%  '
$ws.Range('D27').Value = "'" + '8.956'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.24%  '
$ws.Range('D28').Value = "'" + '18.21'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.01%  '
$ws.Range('D29').Value = "'" + '1.345'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.53%  '
$ws.Range('D30').Value = "'" + '1.493'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.73%  '
$ws.Range('D31').Value = "'" + '4.400'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.56%  '
$ws.Range('D32').Value = "'" + '4.103'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.00%  '
$ws.Range('D33').Value = "'" + '0.05194'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.74%  '
$ws.Range('D34').Value = "'" + '1.932'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.75%  '
$ws.Range('D35').Value = "'" + '1.185'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.52%  '
$ws.Range('D36').Value = "'" + '0.7283'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.48%  '
$ws.Range('D37').Value = "'" + '2.680'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('D38').Value = "'" + '0.01855'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.54%  '
$ws.Range('D39').Value = "'" + '2.699'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.34%  '
$ws.Range('D40').Value = '1.174.68'
$ws.Range('E40').Value = '  +0.40%  '
$ws.Range('D41').Value = "'" + '0.9035'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.90%  '
$ws.Range('E42').Value = '  +1.46%  '
$ws.Range('D43').Value = "'" + '72.42'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.59%  '
$ws.Range('D44').Value = "'" + '0.9999'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('D45').Value = "'" + '101.85'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.51%  '
$ws.Range('D46').Value = "'" + '0.5284'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.22%  '
$ws.Range('D47').Value = '2.008.85'
$ws.Range('E47').Value = '  -0.52%  '
$ws.Range('D48').Value = "'" + '1.780'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.15%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = "'" + '0.00000000120'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.27%  '
$ws.Range('B50').Value = 'SynthetixNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D50').Value = "'" + '2.893'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.83%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = "'" + '9.296'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.93%  '
